# The single "Localización" column (D) holding a combined "41.5N35.99W"
# string is split into two separate numeric columns: "Latitud" and
# "Longitud". This pushes the existing "Tipo" column from E to F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at E, shifting the old E ("Tipo") to F and
# leaving D ("Localización" header / "41.5N35.99W" value) untouched for now.
$ws.Range("E1").EntireColumn.Insert()

# Replace the old "Localización" header with "Latitud" and give the new
# column the "Longitud" header.
$ws.Range("D1").Value = "Latitud"
$ws.Range("E1").Value = "Longitud"

# Replace the combined coordinate text with the split numeric values.
$ws.Range("D2").Value = 12.6987
$ws.Range("E2").Value = 15.268

# Match the author's resulting selection.
$ws.Range("D2").Select()
